$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text, matching the source data's text-based
# formatting (values like '36.452.55' or '242.89' are not valid numbers).
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '36.396.76'
$ws.Range('E2').Value = '  -0.37%  '
$ws.Range('D3').Value = '1.941.50'
$ws.Range('E3').Value = '  -2.21%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '242.67'
$ws.Range('E5').Value = '  -0.83%  '
$ws.Range('D6').Value = '0.611'
$ws.Range('E6').Value = '  -2.81%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '57.22'
$ws.Range('E8').Value = '  -3.63%  '
$ws.Range('D9').Value = '0.362'
$ws.Range('E9').Value = '  -3.66%  '
$ws.Range('D10').Value = '0.0851'
$ws.Range('E10').Value = '  +3.09%  '
$ws.Range('D11').Value = '0.103'
$ws.Range('E11').Value = '  -0.60%  '
$ws.Range('D12').Value = '2.226.58'
$ws.Range('E12').Value = '  -2.17%  '
$ws.Range('D13').Value = '0.816'
$ws.Range('E13').Value = '  -5.82%  '
$ws.Range('D14').Value = '21.14'
$ws.Range('E14').Value = '  -10.86%  '
$ws.Range('D15').Value = '13.48'
$ws.Range('E15').Value = '  -4.02%  '
$ws.Range('D16').Value = '5.19'
$ws.Range('E16').Value = '  -5.30%  '
$ws.Range('D17').Value = '1.953.13'
$ws.Range('E17').Value = '  -1.52%  '
$ws.Range('D18').Value = '36.354.25'
$ws.Range('E18').Value = '  -0.12%  '
$ws.Range('D19').Value = '0.0₃0873'
$ws.Range('E19').Value = '  +0.94%  '
$ws.Range('D20').Value = '69.36'
$ws.Range('E20').Value = '  -1.58%  '
$ws.Range('D21').Value = '228.66'
$ws.Range('E21').Value = '  -2.34%  '
$ws.Range('D22').Value = '5.01'
$ws.Range('E22').Value = '  -6.26%  '
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('D24').Value = '2.39'
$ws.Range('E24').Value = '  -8.10%  '
$ws.Range('D25').Value = '2.28'
$ws.Range('E25').Value = '  -1.29%  '
$ws.Range('D26').Value = '9.27'
$ws.Range('E26').Value = '  -8.29%  '
$ws.Range('D27').Value = '160.93'
$ws.Range('E27').Value = '  -0.83%  '
$ws.Range('D28').Value = '0.131'
$ws.Range('E28').Value = '  -0.31%  '
$ws.Range('D29').Value = '19.26'
$ws.Range('E29').Value = '  -3.10%  '
$ws.Range('D30').Value = '0.118'
$ws.Range('E30').Value = '  -2.15%  '
$ws.Range('D31').Value = '1.13'
$ws.Range('E31').Value = '  -5.70%  '
$ws.Range('D32').Value = '4.62'
$ws.Range('E32').Value = '  -6.21%  '
$ws.Range('D33').Value = '0.0635'
$ws.Range('E33').Value = '  +0.75%  '
$ws.Range('D34').Value = '4.23'
$ws.Range('E34').Value = '  -4.43%  '
$ws.Range('D35').Value = '6.13'
$ws.Range('E35').Value = '  -3.18%  '
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  +0.16%  '
$ws.Range('D37').Value = '1.79'
$ws.Range('E37').Value = '  +0.68%  '
$ws.Range('D38').Value = '2.15'
$ws.Range('E38').Value = '  -5.43%  '
$ws.Range('D39').Value = '3.06'
$ws.Range('E39').Value = '  +0.66%  '
$ws.Range('D40').Value = '0.0978'
$ws.Range('E40').Value = '  +1.77%  '
$ws.Range('D41').Value = '2.86'
$ws.Range('E41').Value = '  -1.39%  '
$ws.Range('D42').Value = '1.16'
$ws.Range('E42').Value = '  -6.51%  '
$ws.Range('D43').Value = '0.0210'
$ws.Range('E43').Value = '  -1.72%  '
$ws.Range('D44').Value = '15.73'
$ws.Range('E44').Value = '  -3.14%  '
$ws.Range('D45').Value = '1.341.27'
$ws.Range('E45').Value = '  -3.07%  '
$ws.Range('D46').Value = '1.02'
$ws.Range('E46').Value = '  -6.67%  '
$ws.Range('D47').Value = '87.20'
$ws.Range('E47').Value = '  -6.14%  '
$ws.Range('D48').Value = '7.14'
$ws.Range('E48').Value = '  -5.23%  '
$ws.Range('D49').Value = '2.83'
$ws.Range('E49').Value = '  -0.85%  '
$ws.Range('D50').Value = '44.34'
$ws.Range('E50').Value = '  -2.67%  '
$ws.Range('D51').Value = '2.119.26'
$ws.Range('E51').Value = '  -2.26%  '
